$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change A1 label from "NON STEMMED" to "STEMMED" (this is the header for the
# first results block, which originally duplicated the "NON STEMMED" label
# used lower down for the non-stemmed proximity-search block in row 12).
$ws.Range("A1").Value = "STEMMED"

# Move the active selection to A4 (was A8 before the edit).
$ws.Range("A4").Select()
